# Worked on temporal resolution
# Extend the "Demand" sheet's time series (column B, commodity EU27.Elec)
# from a single timestep (t=1) to a full 12-timestep profile (t=1..12),
# and correct the demand value used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Row 3 (t=1): replace the old demand figure with the corrected value.
$ws.Range("B3").Value = 345838542

# Rows 4-14 (t=2..12): new timesteps, same corrected demand value.
for ($i = 4; $i -le 14; $i++) {
    $t = $i - 2
    $ws.Cells.Item($i, 1).Value = $t
    $ws.Cells.Item($i, 2).Value = 345838542
}

# Column B only needs to be wide enough to show the new values.
$ws.Columns.Item(2).ColumnWidth = 9.140625

# Make "Demand" the active sheet/tab with the new data selected,
# matching where the author left off after extending the series.
$ws.Range("B3:B14").Select()
